# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" on all sheets
# - Narrow the "Status" column(s) now that the text is shorter

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update the status text wherever it appears
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Narrow the Status columns to match the new, shorter content
$overview.Range("E1").EntireColumn.ColumnWidth = 12.5
$overview.Range("F1").EntireColumn.ColumnWidth = 12.5
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
